$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leetcode刷题")

# Row 147: 122. 买卖股票的最佳时机 II
$ws.Range("C147").Value = "122. 买卖股票的最佳时机 II"
$ws.Range("D147").Value = "122. 买卖股票的最佳时机 II"
$ws.Range("E147").Value = 9.24
$ws.Range("G147").Value = "js:`n代码好写,就是思路,没有抽象得很干净`n本题中理解利润拆分是关键点！ 不要整块的去看，而是把整体利润拆为每天的利润。`n一旦想到这里了，很自然就会想到贪心了，即：只收集每天的正利润，最后稳稳的就是最大利润了。"
$ws.Rows.Item(147).RowHeight = 60

# Row 148: 55. 跳跃游戏
$ws.Range("C148").Value = "55. 跳跃游戏"
$ws.Range("D148").Value = "55. 跳跃游戏"
$ws.Range("E148").Value = 9.25
$ws.Range("G148").Value = "js:`n有思路,写了半天模拟出来,发现贪心咋不行??原来是自己的贪心思路错了...`n贪心算法局部最优解：每次取最大跳跃步数（取最大覆盖范围），整体最优解：最后得到整体最大覆盖范围，看是否能到终点。"
$ws.Rows.Item(148).RowHeight = 36

# Row 149: 45. 跳跃游戏 II
$ws.Range("C149").Value = "45. 跳跃游戏 II"
$ws.Range("D149").Value = "45. 跳跃游戏 II"
$ws.Range("E149").Value = 9.25
$ws.Range("G149").Value = "js:`n没有思路, 思路巧妙, 代码简洁, 抄的 :  以最小的步数增加最大的覆盖范围，直到覆盖范围覆盖了终点，这个范围内最少步数一定可以跳到，不用管具体是怎么跳的，不纠结于一步究竟跳一个单位还是两个单位"
$ws.Rows.Item(149).RowHeight = 36

# Row 150: 1005. K 次取反后最大化的数组和
$ws.Range("C150").Value = "1005. K 次取反后最大化的数组和"
$ws.Range("D150").Value = "1005. K 次取反后最大化的数组和"
$ws.Range("E150").Value = 9.25
$ws.Range("G150").Value = "js:`n方法一: 自己写的. 天哪,好久没有自己oc一个题了,开心!`n方法二: 优化`n做完这个题不是很能感觉用到了贪心,实际上用了两次,要梳理思路,贪心的思路要有"
$ws.Rows.Item(150).RowHeight = 60

# Update view: scroll/selection to show the new rows
$ws.Range("G150").Select()
